$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "CasiPlus"
$ws.Range("C10").Value = "Skamol"
$ws.Range("D10").Value = "25, 30, 50"
$ws.Range("A10").Value = "903, 125, 705"
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 8

$ws.Range("C20").Select()
